$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A74").Value = 900
$ws.Range("B74").Value = 520
$ws.Range("C74").Value = 40
$ws.Range("D74").Value = 102.052852566979
